# Insert a new data row at row 431 (pushing the existing rows 431-514 down
# to 432-515) and populate it with the new observation, matching the
# weekly data refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(431).Insert()

$ws.Cells.Item(431, 1).Value = 4
$ws.Cells.Item(431, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(431, 3).Value = "Los Lagos"
$ws.Cells.Item(431, 4).Value = 45211
$ws.Cells.Item(431, 5).Value = 10
$ws.Cells.Item(431, 6).Value = 100112043
$ws.Cells.Item(431, 7).Value = "Pepino ensalada"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 200
$ws.Cells.Item(431, 11).Value = 20000
$ws.Cells.Item(431, 12).Value = 20000
$ws.Cells.Item(431, 13).Value = 20000
$ws.Cells.Item(431, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(431, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(431, 16).Value = 333
$ws.Cells.Item(431, 17).Value = 60
$ws.Cells.Item(431, 18).Value = "Hortaliza"
